# Update HoverNet and stain separation descriptions (tools_summary.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2 (HistoReg row, "Notes"): append a period before the trailing space.
$ws.Range("F2").Value = "I encountered some difficulties building ITK as described in installation instructions. Issues were resolved on Ubuntu 20.04 with local installation of ITK 4.13.2 using  apt package. "

# F4 (HoverNet row, "Notes"): previously empty, now has a comparison note.
$ws.Range("F4").Value = "Empirically, appears to be a good middle ground between StarDist and CellPose in that segmentations are not confined to star-convex polygons (StarDist) but do not overshoot the nucleus (CellPose)."

# F5 (StarDist row, "Notes"): extend with a caveat about non-rounded objects.
$ws.Range("F5").Value = "All detected objects will take rounded, star-convex polygon structure. Tool is therefore valid on rounded nuclei but not so much on other objects that may take elongated or ‘sharp’ shapes."

# F10 (Staintools stain-extraction row, "Notes"): reworded.
$ws.Range("F10").Value = "Separation into two stains i.e. H&E into H and E or H-DAB into H and DAB."

# F11 (Geijs et al. AutomaticColorUnmixing row, "Notes"): reworded to three stains.
$ws.Range("F11").Value = "Separation into three stains i.e. H, E (residual) and DAB; GitHub repo is a public fork of an internal one developed by Geijs et al. Repo editing to increase dependency on openly accessible, stable libraries (scikit-image etc.) ongoing – contributions welcome"

# F12 (Ruifrok & Johnston / skimage rgb2hed row, "Notes"): reworded to three stains.
$ws.Range("F12").Value = "Separation into three stains i.e. H, E (residual) and DAB."

# Row heights grow to fit the newly lengthened "Notes" text.
$ws.Rows.Item(4).RowHeight = 122.5
$ws.Rows.Item(5).RowHeight = 103.65
$ws.Rows.Item(11).RowHeight = 128

# Match the author's final selection/scroll position in the saved file.
[void]$ws.Range("F11").Select()
